# Recalculating true V_A and lost_Va
#
# 1. Rename the existing (only) sheet and add a new second sheet.
# 2. Populate the new sheet with a small summary table.
# 3. Update the job-ID note on the first sheet (Set_15_B / I21).
# 4. Leave the new sheet as the active / selected tab, matching the
#    workbookView's activeTab + the new sheet's tabSelected state.

$wb = $excel.ActiveWorkbook

# --- Sheet1: rename + tweak a log entry -----------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sims and main analysis"

$ws1.Range("I21").Value = "job_262387 (ac3) for 1-65; 284708 (qm) for 66-100 [the ac3 job array stopped because of maintenance on ac3]"

# Selection on sheet1 moves to D21 (no longer the active/tabSelected sheet).
$ws1.Range("D21").Select() | Out-Null

# --- Sheet2: new sheet with recalculation summary table --------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Rcalculate true V_A & V_a_left"

# Column widths (character units) to match the source layout as closely
# as the host's width quantization allows.
$ws2.Columns.Item(2).ColumnWidth = 9.436197916666666
$ws2.Columns.Item(3).ColumnWidth = 9.436197916666666
$ws2.Columns.Item(4).ColumnWidth = 15.256510416666666

# Header row (bold, reuses the workbook's existing bold style).
$ws2.Range("A1:F1").Font.Bold = $true
$ws2.Range("A1").Value = "Set(s)"
$ws2.Range("B1").Value = "Description"
$ws2.Range("C1").Value = "nsims"
$ws2.Range("D1").Value = "Start Date"
$ws2.Range("E1").Value = "End Date"
$ws2.Range("F1").Value = "Job ID"

# Data row.
$ws2.Range("A2").Value = "Set_9, Set_15a, Set_15_b"
$ws2.Range("B2").Value = "Standard set for sims with burnin and sims with larger scales"
$ws2.Range("C2").Value = 300
$ws2.Range("D2").NumberFormat = "d-mmm-yy"
$ws2.Range("D2").Value = "3/7/2025"

$ws2.Range("D3").Select() | Out-Null
